# Chapter 6 edit:
#   - insert a new "1.1 subtitle" Heading2 paragraph right after the
#     chapter Heading1 ("Human-Centricity: Design Constraints and
#     Opportunities")
#   - wrap that new paragraph plus the following FirstParagraph
#     paragraph in a new "subtitle" bookmark (nested inside the
#     existing "chapter-6" bookmark)
#   - change the placeholder body text "[Target x words]" to "test"

$d = $word.ActiveDocument

# --- locate the chapter heading paragraph -------------------------------
$heading1 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "*Human-Centricity*") {
        $heading1 = $candidate
    }
}

# --- insert the new "1.1`tsubtitle" paragraph right after it -----------
$headingIndex = $heading1.Index
$insertionPoint = $heading1.Range
$insertionPoint.Collapse(0)                 # wdCollapseEnd
$insertionPoint.InsertAfter("1.1`tsubtitle`r")

$subtitlePara = $d.Paragraphs.Item($headingIndex + 1)
$subtitlePara.Style = "Heading2"

# Give the leading "1.1" its own run with the SectionNumber character style,
# matching the way the chapter number itself is styled.
$numberRange = $d.Range($subtitlePara.Range.Start, $subtitlePara.Range.Start + 3)
$numberRange.Style = "SectionNumber"

# --- swap the placeholder paragraph text --------------------------------
$d.Content.Find.Execute("[Target x words]", $false, $false, $false, $false,
                         $false, $true, 1, $false, "test", 2)

# --- wrap the new subtitle paragraph + the body paragraph in a bookmark -
$subtitlePara = $d.Paragraphs.Item($headingIndex + 1)   # "1.1`tsubtitle"
$bodyPara = $d.Paragraphs.Item($headingIndex + 2)       # "test"
$bookmarkRange = $d.Range($subtitlePara.Range.Start, $bodyPara.Range.End)
$d.Bookmarks.Add("subtitle", $bookmarkRange)

Write-Host "chapter-6 subtitle section inserted"
